$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.922.09"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.745.05"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'246.97"
$ws.Range("E5").Value = "  +3.86%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.5023"
$ws.Range("E7").Value = "  -4.80%  "
$ws.Range("D8").Value = "'0.2732"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "'0.06176"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "1.746.38"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "'0.07255"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'0.6522"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'15.12"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").Value = "'4.633"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "'77.61"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'0.9997"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "25.952.81"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'0.000006806"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "1.969.25"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "'4.342"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'8.666"
$ws.Range("D24").Value = "'5.395"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").Value = "'137.29"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("D26").Value = "'1.501"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("D27").Value = "'15.21"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").Value = "'105.54"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'3.911"
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("D31").Value = "'0.08228"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "'3.636"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'0.04667"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "'0.9927"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").Value = "'0.6184"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").Value = "'0.01603"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").Value = "'1.914"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("D40").Value = "'1.000"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'99.74"
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("D42").Value = "'0.3886"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").Value = "'0.7563"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "'5.000"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "'0.1144"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'6.288"
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("D47").Value = "'55.44"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").Value = "'0.05246"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").Value = "'30.61"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("D50").Value = "'7.551"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "'0.3426"
$ws.Range("E51").Value = "  -1.73%  "
